$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 0.2853347213291343
$ws.Range("J3").Value = 0.5652247219704619
$ws.Range("K3").Value = 0.47170424782565
$ws.Range("L3").Value = 2.699537866382877
